$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently lists 8 molecules in A2:A9 (A1 is the "Molecule" header).
# Keep only the "CEFTRIAXONE 1GM" row (originally A3) and drop the rest, so the
# sheet ends up with just the header in A1 and "CEFTRIAXONE 1GM" in A2.

# Remove A2 (MEROPENEM 1GM); this shifts CEFTRIAXONE 1GM up from A3 to A2.
$ws.Rows(2).Delete()

# Remove the remaining unwanted rows (now at A3:A8).
$ws.Range("A3:A8").EntireRow.Delete()

# Match the saved selection state.
$ws.Range("A2").Select() | Out-Null
